# Update loading_percent values for the "case with 380 kV" run
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2  = @{ B=30.24164973136757; C=20.12138428807692; D=11.32638789073647; E=9.421041717961881;  G=3.896815679667412; I=60.84471998896181; J=6.339821035386126; L=15.29519114330293; M=25.29146520893755 }
    3  = @{ B=30.30763363416867; C=19.84288164637158; D=11.3448539259029;  E=9.410500020831561;  G=3.902670446694226; I=59.79935659241307; J=6.340085918273412; L=15.31002555099842; M=25.36201023936074 }
    4  = @{ B=30.36026486125808; C=19.67649898592092; D=11.35774628189282; E=9.403876125696749;  G=3.906441012234904; I=59.14873965907228; J=6.349443564193116; L=15.32149768157409; M=25.41183662861178 }
    5  = @{ B=30.38474131699237; C=19.60993535024922; D=11.36339055464187; E=9.401138328781796;  G=3.908021957760292; I=58.88159970577013; J=6.35356870525961;  L=15.32676618510613; M=25.4337729921229  }
    6  = @{ B=30.38898798712928; C=19.59895951764748; D=11.36435136087871; E=9.400681380300806;  G=3.908287161344809; I=58.83712578707731; J=6.354272491515522; L=15.32767683591035; M=25.43751388981473 }
    7  = @{ B=30.36058272072404; C=19.67559617039201; D=11.35782082163016; E=9.403839359343348;  G=3.906462153308825; I=59.14514478782116; J=6.34949793577874;  L=15.3215663325208;  M=25.41212587094617 }
    8  = @{ B=30.26187572848103; C=20.02444499926342; D=11.33243226246247; E=9.417438045341438;  G=3.898798069410044; I=60.48621192828474; J=6.331168552873715; L=15.29981495512768; M=25.31443434550252 }
    9  = @{ B=30.16518868425662; C=20.74167903657744; D=11.29498898250236; E=9.442927082315194;  G=3.885152769407968; I=63.03827148313088; J=6.42987062885391;  L=15.27595283614129; M=25.17477273901355 }
    10 = @{ B=30.1541061857906;  C=21.28403330501336; D=11.2750215958141;  E=9.460967680900653;  G=3.87595682402209;  I=64.85525080256589; J=6.503346515624683; L=15.26992581675879; M=25.10413519282277 }
    11 = @{ B=30.16223589815105; C=21.5330794485449;  D=11.26757900891726; E=9.469032843464781;  G=3.8719502970093;   I=65.66713815425838; J=6.535971570592015; L=15.26969051762521; M=25.07900739318517 }
    12 = @{ B=30.16721775209826; C=21.62763603970193; D=11.2649968434332;  E=9.472067103623631;  G=3.870458301477048; I=65.97230368098647; J=6.548210004689064; L=15.26996228884329; M=25.07050434228304 }
    13 = @{ B=30.16606004424374; C=21.60726190322717; D=11.26554245099994; E=9.471414495534782;  G=3.87077851274639;  I=65.90668477321479; J=6.545579412158616; L=15.26988770203392; M=25.07229052907309 }
    14 = @{ B=30.16260756720717; C=21.54085419012142; D=11.26736183783077; E=9.469282868528014;  G=3.871827046117284; I=65.69229091266192; J=6.536980763677549; L=15.26970564175243; M=25.07828752702282 }
    15 = @{ B=30.16074093315579; C=21.50020733958214; D=11.26850702959866; E=9.467974617000602;  G=3.872472576690162; I=65.56066685732458; J=6.531698711850977; L=15.26964113257929; M=25.08209283341008 }
    16 = @{ B=30.1538407706003;  C=21.26779695475187; D=11.2755410231908;  E=9.460437798513198;  G=3.876222197291185; I=64.80188223811491; J=6.50119821752108;  L=15.26999167100938; M=25.10591875016478 }
    17 = @{ B=30.15298847626847; C=21.12575470060897; D=11.28027655504318; E=9.455778305499235;  G=3.878567577432532; I=64.33251271795396; J=6.482281661590948; L=15.2708490273762;  M=25.12233306170828 }
    18 = @{ B=30.15373756740346; C=21.04427982254739; D=11.28315472246214; E=9.453084929107089;  G=3.879933228733357; I=64.06117298975448; J=6.471326164724125; L=15.27157805235259; M=25.13243326292158 }
    19 = @{ B=30.15420373168216; C=21.01673496734498; D=11.28415573325882; E=9.452170687989582;  G=3.880398481486627; I=63.96907189875364; J=6.467603938033135; L=15.27186538744533; M=25.1359660679323  }
    20 = @{ B=30.15295086684522; C=21.14085278088829; D=11.27975646634695; E=9.456275695376497;  G=3.87831618607818;  I=64.38262103381722; J=6.48430315091451;  L=15.27073334233155; M=25.12051747696837 }
    21 = @{ B=30.16356992420807; C=21.56035368843195; D=11.26682102804247; E=9.469909513376727;  G=3.871518384639719; I=65.75532673524724; J=6.539509550275263; L=15.2697493201887;   M=25.07649855216199 }
    22 = @{ B=30.18160687952318; C=21.83592984121298; D=11.25974371346479; E=9.478704594876646;  G=3.867222344368484; I=66.63911549753031; J=6.574913743305403; L=15.27120978652146; M=25.05363216778881 }
    23 = @{ B=30.17096232094169; C=21.68874889432006; D=11.26339495996743; E=9.474020880360499;  G=3.869501872956272; I=66.16869642604065; J=6.55608007008675;  L=15.27023770894065; M=25.06529476310197 }
    24 = @{ B=30.15296401131782; C=21.13402634943456; D=11.27999111365997; E=9.456050870840228;  G=3.878429786327994; I=64.35997170141654; J=6.483389485061257; L=15.27078490811681; M=25.12133623687799 }
    25 = @{ B=30.18087420649012; C=20.54459273843417; D=11.3037950495293;  E=9.436154454835362;  G=3.888697512751889; I=62.35731884439468; J=6.402110161014559; L=15.28039094065232; M=25.20696283408759 }
}

foreach ($rowNum in $data.Keys) {
    $rowValues = $data[$rowNum]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$rowNum").Value = $rowValues[$col]
    }
}

$wb.Save()
